$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new merged header block D2:H2, styled like A2:C2 but with its own accent fill ---
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("D2:H2").PasteSpecial(-4122) | Out-Null
$ws.Range("D2:H2").Interior.Color = 11919046
$ws.Range("D2:H2").Merge() | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("D2").Value = "Comandos Consola"

# --- Rows 3-7: rebuild the project reference / console commands table ---
# Row 3 - DIAGRAMA
$ws.Range("A3").Value = "DIAGRAMA"
$ws.Range("B3").Value = "UML"
$ws.Range("C3").Value = "draw.io"

# Row 4 - ARCHIVO SQL
$ws.Range("A4").Value = "ARCHIVO SQL"
$ws.Range("B4").Value = "BDSQLServer"
$ws.Range("C4").Value = "SQL Server"

# Row 5 - PROYECTO
$ws.Range("A5").Value = "PROYECTO"
$ws.Range("B5").Value = "WSVenta"
$ws.Range("C5").Value = "AP.NET Core Web API"
$ws.Range("D5").Value = 'Scaffold-DBContext "Server=OFITE-GRUDE8\SQLEXPRESS;Database=VentaReal;Trusted_Connection=True;" Microsoft.EntityFrameworkCore.SqlServer -OutputDir Models'
$ws.Range("E5").Value = "Conexión BD con autentificacion de windows"

# Row 6 - nugets / second scaffold command
$ws.Range("B6").Value = "nugets"
$ws.Range("C6").Value = "Microsoft.EntityFrameworkCore.SqlServer 3.1.21"
$ws.Range("D6").Value = 'Scaffold-DBContext "Server=OFITE-GRUDE8\SQLEXPRESS;Database=VentaReal;Trusted_Connection=False;user=sa;Password=abc123;" Microsoft.EntityFrameworkCore.SqlServer -OutputDir Models'
$ws.Range("E6").Value = "Conexión BD sin autentificacion de windows"

# Row 7 - second nuget package
$ws.Range("C7").Value = "Microsoft.EntityFrameworkCore.Tools 3.1.21"

# --- Column widths so the new text is fully visible ---
$ws.Columns.Item(1).ColumnWidth = 12.0
$ws.Columns.Item(3).ColumnWidth = 43.83333333333333
$ws.Columns.Item(4).ColumnWidth = 17.166666666666668
$ws.Columns.Item(5).ColumnWidth = 40.666666666666664

# --- Final cursor position, matching the author's last selection ---
$ws.Range("E9").Select() | Out-Null
